$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 8 with data for scenario ID_0010
$ws.Range("A8").Value = "ID_0010"
$ws.Range("B8").Value = "André Automatizador"
$ws.Range("C8").Value = "sem email"
$ws.Range("D8").Value = "automacaoteste"

# Update selection to the newly added cell B8 (as in the edited file)
[void]$ws.Range("B8").Select()
